$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237 - this shifts rows 237..352 down to 238..353
# (Excel copies the formatting of the row above down into the new row, which
# matches the target: D column keeps its date number format.)
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new weekly record.
$ws.Range("A237").Value = 4
$ws.Range("B237").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C237").Value = "Los Lagos"
$ws.Range("D237").Value = 45016
$ws.Range("E237").Value = 10
$ws.Range("F237").Value = "Fruta"
$ws.Range("G237").Value = 100108
$ws.Range("H237").Value = "Tropicales y subtropicales"
$ws.Range("I237").Value = 100108002
$ws.Range("J237").Value = "Mango"
$ws.Range("K237").Value = "Sin especificar"
$ws.Range("L237").Value = "Primera"
$ws.Range("M237").Value = 200
$ws.Range("N237").Value = 8000
$ws.Range("O237").Value = 8500
$ws.Range("P237").Value = 8250
$ws.Range("Q237").Value = "$/bandeja 4 kilos"
$ws.Range("R237").Value = "Perú"
$ws.Range("S237").Value = 2062
$ws.Range("T237").Value = 4
